$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row after "Contact" (row 10) for "Jurisdiction", copying the
# formatting of the row above so the new row matches the existing style.
$ws.Rows.Item(11).Insert(-4121, 0)
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Update Version value (now row 3)
$ws.Cells.Item(3, 2).Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (now row 8)
$ws.Cells.Item(8, 2).Value = "2025-10-29T22:15:57+01:00"
